# Atualização de bases das ligas, do dia: 11-04-2024 às 23:56
#
# For a set of match-row pairs, the entire data payload (columns B..AC —
# id, Div, Div Original Name, HomeTeam, AwayTeam, score, odds, ...) was
# swapped between the two rows of each pair, while the leading sequence
# index in column A stayed put on its own row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(19, 20),
    @(38, 39),
    @(50, 51),
    @(77, 78),
    @(84, 85),
    @(90, 91),
    @(140, 141),
    @(188, 189),
    @(195, 196),
    @(212, 213),
    @(214, 215)
)

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}
